$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 94 and row 95 (columns A and B), keep formatting
$ws.Range("A94").Value = ""
$ws.Range("B94").Value = ""
$ws.Range("A95").Value = ""
$ws.Range("B95").Value = ""

# Add new "erledigte Tasks" notes in column D
# (order matters for shared-string table index assignment, so add them
# in the same order as the original authoring tool did)
$ws.Range("C103").Value = "Mikula"
$ws.Range("D103").Value = "Gui Admin Add Product fertig"
$ws.Range("D100").Value = "GUI admin restock article fertig"
$ws.Range("D99").Value = "Gui User Order fertig"

# Update the active view (scroll position / selection)
$ws.Range("F97").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 76
